$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.46"
$ws.Range("E2").Value = "'0.16%"

$ws.Range("D3").Value = "'37.07"
$ws.Range("E3").Value = "'6.48%"

$ws.Range("D4").Value = "'5.006"
$ws.Range("E4").Value = "'-3.15%"

$ws.Range("D5").Value = "'0.07876"
$ws.Range("E5").Value = "'1.02%"

$ws.Range("D6").Value = "'2.206"
$ws.Range("E6").Value = "'-4.27%"

$ws.Range("D7").Value = "'8.013"
$ws.Range("E7").Value = "'-0.53%"

$ws.Range("D8").Value = "'4.020"
$ws.Range("E8").Value = "'0.78%"

$ws.Range("D9").Value = "'0.9205"
$ws.Range("E9").Value = "'-0.45%"

$ws.Range("D10").Value = "'0.09642"
$ws.Range("E10").Value = "'-4.10%"

$ws.Range("D11").Value = "'0.1890"
$ws.Range("E11").Value = "'3.49%"

$ws.Range("D12").Value = "'0.08582"
$ws.Range("E12").Value = "'0.96%"

$ws.Range("D13").Value = "'0.03686"
$ws.Range("E13").Value = "'8.74%"

$ws.Range("D14").Value = "'0.09980"
$ws.Range("E14").Value = "'0.66%"

$ws.Range("D15").Value = "'0.001495"
$ws.Range("E15").Value = "'0.02%"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04554"
$ws.Range("E16").Value = "'-2.10%"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005689"
$ws.Range("E17").Value = "'-2.04%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.469"
$ws.Range("E18").Value = "'-0.07%"

$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.249"
$ws.Range("E19").Value = "'6.99%"

$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3414"
$ws.Range("E20").Value = "'-0.05%"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1317"
$ws.Range("E21").Value = "'-0.72%"

$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'4.750"
$ws.Range("E22").Value = "'4.29%"

$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.2200"
$ws.Range("E23").Value = "'-3.22%"

$ws.Range("D24").Value = "'0.001234"
$ws.Range("E24").Value = "'1.22%"

$ws.Range("D25").Value = "'0.004470"
$ws.Range("E25").Value = "'3.02%"

$ws.Range("D26").Value = "'0.0001399"
$ws.Range("E26").Value = "'7.43%"

$ws.Range("E27").Value = "'39.61%"

$ws.Range("D39").Value = "'0.01841"
$ws.Range("E39").Value = "'5.53%"

$ws.Range("E40").Value = "'0.12%"

$ws.Range("D41").Value = "'0.008124"
$ws.Range("E41").Value = "'5.73%"

$ws.Range("D42").Value = "'0.1398"
$ws.Range("E42").Value = "'-0.97%"

$ws.Range("D43").Value = "'0.007541"
$ws.Range("E43").Value = "'-1.42%"

$ws.Range("D44").Value = "'0.002228"
$ws.Range("E44").Value = "'-2.86%"

$ws.Range("D45").Value = "'0.01050"
$ws.Range("E45").Value = "'5.32%"

$ws.Range("D46").Value = "'0.00006269"
$ws.Range("E46").Value = "'3.33%"

$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.23%"

$ws.Range("D48").Value = "'0.0005799"
$ws.Range("E48").Value = "'-0.02%"

$ws.Range("D49").Value = "'29.85"
$ws.Range("E49").Value = "'415.00%"

$ws.Range("E50").Value = "'-36.12%"

$ws.Range("D51").Value = "'0.00002098"
$ws.Range("E51").Value = "'-0.23%"
